# penambahan temp table untuk api hotel
$wb = $excel.ActiveWorkbook

# --- Add the new "Temp" worksheet after the last existing sheet (Log) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tempSheet = $wb.Worksheets.Add($null, $lastSheet)
$tempSheet.Name = "Temp"

# Write values in the same order they were originally entered so the
# shared-string table comes out in the same sequence (temp001, Kode table,
# nama Table, Temp Header Hotel API).
$tempSheet.Range("A2").Value = "temp001"
$tempSheet.Range("A1").Value = "Kode table"
$tempSheet.Range("B1").Value = "nama Table"
$tempSheet.Range("B2").Value = "Temp Header Hotel API"

# Match the column widths used for the header labels on the new sheet.
$tempSheet.Columns.Item(1).ColumnWidth = 9.736979166666666
$tempSheet.Columns.Item(2).ColumnWidth = 10.307291666666666

# --- Update the "Log" sheet header row to reuse the same new labels ---
$logSheet = $wb.Worksheets.Item("Log")
$logSheet.Range("A1").Value = "Kode table"
$logSheet.Range("B1").Value = "nama Table"
$logSheet.Range("A1:B1").Select() | Out-Null

# --- Leave the new "Temp" sheet selected/active, matching the commit ---
$tempSheet.Range("B3").Select() | Out-Null
$tempSheet.Activate() | Out-Null
